# Refined metadata to be additional tab
#
# 1) Update the "time_taken" (column F) timestamps on the "data" sheet to
#    reflect the re-run query time.
# 2) Add a new "metadata" worksheet (after "data") describing the panel
#    query that produced the "data" sheet.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- 1) refresh the time_taken column on "data" ------------------------
$data.Range("F2").Value  = "2021-10-05 14:22:07.458854"
$data.Range("F3").Value  = "2021-10-05 14:22:07.458865"
$data.Range("F4").Value  = "2021-10-05 14:22:07.458868"
$data.Range("F5").Value  = "2021-10-05 14:22:07.458871"
$data.Range("F6").Value  = "2021-10-05 14:22:07.458874"
$data.Range("F7").Value  = "2021-10-05 14:22:07.458877"
$data.Range("F8").Value  = "2021-10-05 14:22:07.458880"
$data.Range("F9").Value  = "2021-10-05 14:22:07.458883"
$data.Range("F10").Value = "2021-10-05 14:22:07.458886"
$data.Range("F11").Value = "2021-10-05 14:22:07.458889"
$data.Range("F12").Value = "2021-10-05 14:22:07.458891"
$data.Range("F13").Value = "2021-10-05 14:22:07.458894"
$data.Range("F14").Value = "2021-10-05 14:22:07.458897"
$data.Range("F15").Value = "2021-10-05 14:22:07.458899"
$data.Range("F16").Value = "2021-10-05 14:22:07.458902"
$data.Range("F17").Value = "2021-10-05 14:22:07.458905"
$data.Range("F18").Value = "2021-10-05 14:22:07.458907"
$data.Range("F19").Value = "2021-10-05 14:22:07.458910"
$data.Range("F20").Value = "2021-10-05 14:22:07.458913"
$data.Range("F21").Value = "2021-10-05 14:22:07.458916"
$data.Range("F22").Value = "2021-10-05 14:22:07.458919"
$data.Range("F23").Value = "2021-10-05 14:22:07.458922"
$data.Range("F24").Value = "2021-10-05 14:22:07.458925"
$data.Range("F25").Value = "2021-10-05 14:22:07.458927"
$data.Range("F26").Value = "2021-10-05 14:22:07.458930"
$data.Range("F27").Value = "2021-10-05 14:22:07.458933"
$data.Range("F28").Value = "2021-10-05 14:22:07.458936"
$data.Range("F29").Value = "2021-10-05 14:22:07.458939"
$data.Range("F30").Value = "2021-10-05 14:22:07.458942"
$data.Range("F31").Value = "2021-10-05 14:22:07.458945"

# --- 2) add the "metadata" sheet, after "data" --------------------------
$meta = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$meta.Name = "metadata"

# Pull the header/row formatting from "data" (same bold/bordered style used
# on its header row and the leading index column) so the new sheet matches
# the workbook's existing look.
$data.Range("B1:F1").Copy($meta.Range("B1:F1"))
$data.Range("F1").Copy($meta.Range("G1"))
$data.Range("A2").Copy($meta.Range("A2"))

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("B2").Value = "Pain syndromes"
$meta.Range("C2").Value = 288
$meta.Range("D2").Formula = "'1.10"
$meta.Range("E2").Value = "2021-08-04T12:55:59.576824Z"
$meta.Range("F2").Value = "2021-10-05 14:22:07.455330"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/288/?format=json"
